# FR-UC crosstable update:
#  - Row 34 ("R-10.2") is retitled to "Q-10.2", absorbing the old Q-10.3 row.
#  - The old row 35 ("Q-10.3") is deleted outright (Excel shifts rows 36-89 up by one,
#    so former FR-11..FR-16 block now sits at rows 35-47, and the trailing blank
#    row 89 disappears).
#  - A batch of new "x" marks is filled in across UC09-UC21 (columns J-V) for
#    requirements FR-06 through FR-16, completing the crosstable.
#  - View state: zoom set to 55% and selection moved to Y52, matching the saved file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the FR-10 sub-requirement row and drop the row below it.
$ws.Range("A34").Value = "Q-10.2"
$ws.Rows(35).Delete()

# Newly completed crosstable marks (rows renumbered post-deletion).
$ws.Range("J21").Value = "x"
$ws.Range("J22").Value = "x"
$ws.Range("J23").Value = "x"

$ws.Range("K24").Value = "x"
$ws.Range("K25").Value = "x"
$ws.Range("K26").Value = "x"

$ws.Range("L27").Value = "x"
$ws.Range("L28").Value = "x"

$ws.Range("M29").Value = "x"
$ws.Range("M30").Value = "x"
$ws.Range("M31").Value = "x"

$ws.Range("N32").Value = "x"
$ws.Range("O32").Value = "x"
$ws.Range("N33").Value = "x"
$ws.Range("O34").Value = "x"

$ws.Range("P35").Value = "x"
$ws.Range("P36").Value = "x"

$ws.Range("Q37").Value = "x"
$ws.Range("Q38").Value = "x"

$ws.Range("R39").Value = "x"
$ws.Range("R40").Value = "x"

$ws.Range("S41").Value = "x"
$ws.Range("S42").Value = "x"

$ws.Range("T43").Value = "x"
$ws.Range("T44").Value = "x"
$ws.Range("V44").Value = "x"
$ws.Range("T45").Value = "x"
$ws.Range("U45").Value = "x"
$ws.Range("T46").Value = "x"
$ws.Range("U46").Value = "x"

$ws.Range("V47").Value = "x"

# Match the saved view state (zoom + last selection).
$excel.ActiveWindow.Zoom = 55
$ws.Range("Y52").Select()
